# Blackjack/EenVariabele.xlsx — add a second ("Perfect" strategy) analysis
# table next to the existing ("Basic" strategy) one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: push the existing table down one row, so row 1 can
#        become a merged "Basic" title row above it (like the new
#        "Perfect" block that will sit in E:G). ------------------------
$ws.Rows("2:2").Insert()

# --- 2. Title row (row 1): two merged, centered headers ----------------
$ws.Range("A1:C1").Merge()
$ws.Range("A1").Value = "Basic"
$ws.Range("A1:C1").HorizontalAlignment = -4108   # xlCenter

$ws.Range("E1:G1").Merge()
$ws.Range("E1").Value = "Perfect"
$ws.Range("E1:G1").HorizontalAlignment = -4108   # xlCenter

# --- 3. "Waarnemingen" column headers (row 2), for both tables ---------
$ws.Range("A2").Value = "Waarnemingen"
$ws.Range("E2").Value = "Waarnemingen"

# --- 4. Raw observations for the "Perfect" strategy (E3:E12) -----------
$ws.Range("E3").Value = 1021.5
$ws.Range("E4").Value = 1016
$ws.Range("E5").Value = 1014
$ws.Range("E6").Value = 1027.5
$ws.Range("E7").Value = 1022.5
$ws.Range("E8").Value = 1017
$ws.Range("E9").Value = 1023
$ws.Range("E10").Value = 1016
$ws.Range("E11").Value = 1019
$ws.Range("E12").Value = 1013.5

# --- 5. Labels + formulas mirroring the existing A:C table --------------
$ws.Range("F3").Value = "Mediaan"
$ws.Range("G3").Formula = "=MEDIAN(E3:E12)"

$ws.Range("F4").Value = "Modus"
$ws.Range("G4").Formula = "=MODE(E3:E12)"

$ws.Range("F5").Value = "Bereik"
$ws.Range("G5").Formula = "=ABS(MAX(E3:E12)-MIN(E3:E12))"

$ws.Range("F6").Value = "Q1"
$ws.Range("G6").Formula = "=QUARTILE(E3:E12,1)"

$ws.Range("F7").Value = "Q2"
$ws.Range("G7").Formula = "=QUARTILE(E3:E12,2)"

$ws.Range("F8").Value = "Q3"
$ws.Range("G8").Formula = "=QUARTILE(E3:E12,3)"

$ws.Range("F9").Value = "Variantie"
$ws.Range("G9").Formula = "=VARP(E3:E12)"

$ws.Range("F10").Value = "standaardafwijking"
$ws.Range("G10").Formula = "=SQRT(G9)"

$ws.Range("F13").Value = "Gemiddelde"
$ws.Range("G13").Formula = "=AVERAGE(E3:E12)"

# --- 6. Column widths for the new block (mirrors A:C) -------------------
$ws.Columns("E").ColumnWidth = $ws.Columns("A").ColumnWidth
$ws.Columns("F").ColumnWidth = $ws.Columns("B").ColumnWidth
$ws.Columns("G").ColumnWidth = $ws.Columns("C").ColumnWidth

# --- 7. View tweaks: scroll right a bit, zoom to 90%, select E13 --------
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Application.ActiveWindow.Zoom = 90
$ws.Range("E13").Select()
